$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Josh Philippe †, Royal Challengers Bangalore): runs 0->1, balls 3->2
# Leading apostrophe forces the numeric-looking value to stay text (matches
# the original cells, which are text-typed numbers), then the style is
# reset to "Normal" so no stray quote-prefix formatting/style index lingers.
$ws.Range("C2").Value = "'1"
$ws.Range("D2").Value = "'2"

# Row 3 (Josh Philippe †, Royal Challengers Bangalore): runs 1->0, balls 2->3
$ws.Range("C3").Value = "'0"
$ws.Range("D3").Value = "'3"

$ws.Range("C2:D3").Style = "Normal"
